$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme
for ($i=9; $i -le 14; $i++) {
  try {
    $c = $cs.Colors($i)
    Write-Output "$i : $($c.RGB)"
  } catch { Write-Output "$i EXC: $_" }
}
